$d = $word.ActiveDocument

# --- Add the three new character styles (appended to styles.xml, same
#     as the diff: GaNStyle, GaNParagraph, GaNLinks, all w:type="character",
#     w:customStyle="1"). WdStyleType 2 == wdStyleTypeCharacter.

$gaNStyle = $d.Styles.Add("GaNStyle", 2)
$gaNStyle.Font.NameAscii = "Calibri"
$gaNStyle.Font.Name = "Calibri"
$gaNStyle.Font.Size = 14

$gaNParagraph = $d.Styles.Add("GaNParagraph", 2)
$gaNParagraph.Font.NameAscii = "Calibri"
$gaNParagraph.Font.Name = "Calibri"
$gaNParagraph.Font.Size = 10

$gaNLinks = $d.Styles.Add("GaNLinks", 2)
$gaNLinks.Font.NameAscii = "Calibri"
$gaNLinks.Font.Name = "Calibri"
$gaNLinks.Font.Bold = $true
$gaNLinks.Font.Color = 8388608
$gaNLinks.Font.Size = 9.5
$gaNLinks.Font.Underline = 1

# --- Apply the new styles to the runs that got a <w:rStyle> in the diff.
#     Each target run already exists (<w:r><w:t>...</w:t></w:r>>); we just
#     need to stamp its character style, which inserts
#     <w:rPr><w:rStyle w:val="..."/></w:rPr> as the run's rPr.

function Set-StyleOnAllMatches($findText, $styleName) {
    $rng = $d.Content
    $rng.Find.ClearFormatting()
    $rng.Find.Text = $findText
    $rng.Find.Forward = $true
    $rng.Find.Wrap = 1
    $rng.Find.MatchCase = $true
    $rng.Find.MatchWholeWord = $false
    $rng.Find.MatchWildcards = $false
    while ($rng.Find.Execute()) {
        $rng.Style = $styleName
        $rng.Collapse(0)
    }
}

Set-StyleOnAllMatches "2022: Daty kampanii używające Gwiazdozbiór Byka: 16-25 stycznia" "GaNStyle"

Set-StyleOnAllMatches "Uczestniczysz w ogólnoświatowym przedsięwzięciu, którego celem jest obserwacja i odnotowanie najsłabszych widocznych gwiazd w celu zmierzenia zanieczyszczenia światłem w danym miejscu. Poprzez zlokalizowanie i obserwację  Gwiazdozbiór Byka na nocnym niebie oraz porównanie go do map nieba ludzie z całego świata będą mogli dowiedzieć się jaki wkład światło emitowane przez ich społeczność wnosi do  zanieczyszczenia światłem. To co dodasz do internetowej bazy danych pomoże udokumentować widoczne nocne niebo." "GaNParagraph"

Set-StyleOnAllMatches " Jenika Hollana, CzechGlobe ((http://amper.ped.muni.cz/jenik/astro/maps/GaNight/2022/)." "GaNLinks"

Write-Output "styles added and applied"
